$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.906
$ws.Range("A9").Value = -21.534
$ws.Range("C9").Value = -10.52
$ws.Range("A18").Value = -21.898
$ws.Range("A20").Value = -20.563
$ws.Range("C23").Value = -12.676
$ws.Range("C24").Value = -12.189
$ws.Range("C26").Value = -12.696
$ws.Range("A27").Value = -21.888
$ws.Range("C34").Value = -11.853
$ws.Range("C35").Value = -12.316
$ws.Range("C48").Value = -11.861
$ws.Range("C52").Value = -11.749
$ws.Range("C66").Value = -11.574
$ws.Range("C67").Value = -10.875
$ws.Range("A69").Value = -21.524
$ws.Range("A76").Value = -20.287
$ws.Range("C80").Value = -12.491
$ws.Range("A82").Value = -22.004
$ws.Range("C99").Value = -11.95
